$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "vsfd"
$ws.Range("E6").Value = "fsd"
$ws.Range("H8").Value = "sdfds"
$ws.Range("N8").Value = "fds"
$ws.Range("J10").Value = 24
$ws.Range("I11").Value = "f"
$ws.Range("D13").Value = "d"
$ws.Range("M13").Value = 425
$ws.Range("E15").Value = "fds"
$ws.Range("L16").Value = 254
$ws.Range("G19").Value = 523

$ws.Range("M13").Select()
